# Updates cryptos list (Coin / Link / Price / Volume(1h)) with refreshed
# quotes pulled on Fri Apr 14 17:54:55 UTC 2023, plus two rows (16/17 and
# 38/39) whose rank swapped places since the last run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes a value while guaranteeing it lands as plain text (Excel would
# otherwise silently coerce numeric-looking strings like "0.9994" or
# "326.88" into a Number cell and lose the original textual formatting).
# Style is reset to "Normal" afterwards so forcing the text number-format
# doesn't leave a stray style on the cell.
function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "30.217.76"
$ws.Range("E2").Value = "  -0.75%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.072.02"
$ws.Range("E3").Value = "  +2.83%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.9994"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5 - BNB
Set-TextValue "D5" "326.88"
$ws.Range("E5").Value = "  +0.58%  "

# Row 6 - USDC
Set-TextValue "D6" "0.9984"
$ws.Range("E6").Value = "  -0.18%  "

# Row 7 - XRP (price unchanged, only volume moved)
$ws.Range("E7").Value = "  +1.35%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.4318"
$ws.Range("E8").Value = "  +4.01%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.08890"
$ws.Range("E9").Value = "  +0.96%  "

# Row 10 - OKB
Set-TextValue "D10" "45.71"
$ws.Range("E10").Value = "  +6.76%  "

# Row 11 - Polygon (price only)
Set-TextValue "D11" "1.152"

# Row 12 - Solana
Set-TextValue "D12" "24.19"
$ws.Range("E12").Value = "  -1.87%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "2.069.71"
$ws.Range("E13").Value = "  +2.86%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.644"
$ws.Range("E14").Value = "  +0.54%  "

# Row 15 - Chainlink
Set-TextValue "D15" "7.641"
$ws.Range("E15").Value = "  +1.92%  "

# Row 16 - was BinanceUSD, now Litecoin (ranking swap with row 17)
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D16" "94.78"
$ws.Range("E16").Value = "  +0.33%  "

# Row 17 - was Litecoin, now BinanceUSD
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D17" "0.9995"
$ws.Range("E17").Value = "  -0.08%  "

# Row 18 - ShibaInu
Set-TextValue "D18" "0.00001118"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06607"
$ws.Range("E19").Value = "  +1.14%  "

# Row 20 - Avalanche
Set-TextValue "D20" "18.72"
$ws.Range("E20").Value = "  -1.66%  "

# Row 21 - Dai (price only)
Set-TextValue "D21" "0.9985"

# Row 22 - Uniswap
Set-TextValue "D22" "6.199"
$ws.Range("E22").Value = "  -0.73%  "

# Row 23 - WrappedBTC
Set-TextValue "D23" "30.250.73"
$ws.Range("E23").Value = "  -0.79%  "

# Row 24 - Cosmos
Set-TextValue "D24" "12.21"
$ws.Range("E24").Value = "  +2.36%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.282"
$ws.Range("E25").Value = "  +2.42%  "

# Row 26 - WrappedliquidstakedEther2.0
Set-TextValue "D26" "2.314.36"
$ws.Range("E26").Value = "  +2.98%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "22.15"
$ws.Range("E27").Value = "  -0.94%  "

# Row 28 - LidoDAOToken
Set-TextValue "D28" "2.525"
$ws.Range("E28").Value = "  +4.11%  "

# Row 29 - Monero
Set-TextValue "D29" "161.50"
$ws.Range("E29").Value = "  -0.92%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "130.76"
$ws.Range("E30").Value = "  -0.59%  "

# Row 31 - ImmutableX (volume only)
$ws.Range("E31").Value = "  +3.91%  "

# Row 32 - Stellar
Set-TextValue "D32" "0.1065"
$ws.Range("E32").Value = "  +1.12%  "

# Row 33 - ARBITRUM
Set-TextValue "D33" "1.623"
$ws.Range("E33").Value = "  +19.39%  "

# Row 34 - Filecoin
Set-TextValue "D34" "6.091"
$ws.Range("E34").Value = "  -0.26%  "

# Row 35 - HuobiToken
Set-TextValue "D35" "3.822"
$ws.Range("E35").Value = "  -0.25%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.02567"
$ws.Range("E36").Value = "  +1.67%  "

# Row 37 - FraxShare
Set-TextValue "D37" "9.655"
$ws.Range("E37").Value = "  +5.72%  "

# Row 38 - was InternetComputer(DFINITY), now Hedera (ranking swap with row 39)
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D38" "0.06613"
$ws.Range("E38").Value = "  -0.10%  "

# Row 39 - was Hedera, now InternetComputer(DFINITY)
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D39" "5.409"
$ws.Range("E39").Value = "  -0.99%  "

# Row 40 - Aptos
Set-TextValue "D40" "12.56"
$ws.Range("E40").Value = "  +1.70%  "

# Row 41 - Algorand
Set-TextValue "D41" "0.2237"
$ws.Range("E41").Value = "  +1.79%  "

# Row 42 - TheSandbox
Set-TextValue "D42" "0.6786"
$ws.Range("E42").Value = "  +1.55%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "1.243"
$ws.Range("E43").Value = "  +0.60%  "

# Row 44 - Frax
Set-TextValue "D44" "0.9979"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "13.88"
$ws.Range("E45").Value = "  +2.07%  "

# Row 46 - Decentraland
Set-TextValue "D46" "0.6327"
$ws.Range("E46").Value = "  +2.21%  "

# Row 47 - NEARProtocol (volume only)
$ws.Range("E47").Value = "  -0.50%  "

# Row 48 - PancakeSwap
Set-TextValue "D48" "3.601"
$ws.Range("E48").Value = "  -1.78%  "

# Row 49 - EOS
Set-TextValue "D49" "1.230"
$ws.Range("E49").Value = "  -3.15%  "

# Row 50 - WEMIXTOKEN
Set-TextValue "D50" "1.189"
$ws.Range("E50").Value = "  +7.05%  "

# Row 51 - Aave
Set-TextValue "D51" "81.10"
$ws.Range("E51").Value = "  -0.55%  "
